# Apply the metrics update described in the diff.
# For every data row (2-26), columns B:Q are overwritten with the same
# new set of metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.6408044419705359,    # B - r2
    -265.1326953808737,    # C - r2_sup
    0.8045844400025129,    # D - r2_test
    0.5069386478537645,    # E - r2_val
    0.7572515948252655,    # F - r2_vt
    0.2132339996864685,    # G - mse
    157.9878643119225,     # H - mse_sup
    0.1194553276864012,    # I - mse_test
    0.08274434155298011,   # J - mse_val
    0.1010998346196906,    # K - mse_vt
    0.2405142646481177,    # L - mape
    0.4617726710043249,    # M - rmse
    0.2163006006629874,    # N - r2_adj
    0.4814312896101858,    # O - rsd
    29.09073025240775,     # P - aic
    44.93611597569436      # Q - bic
)

for ($row = 2; $row -le 26; $row++) {
    for ($i = 0; $i -lt $newValues.Length; $i++) {
        # Column B is index 2 in Cells.Item(row, col)
        $col = $i + 2
        $ws.Cells.Item($row, $col).Value = $newValues[$i]
    }
}
